$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.636.87'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.827.79'
$ws.Range("E3").Value = '  +1.85%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.20'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4680'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3596'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07141'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9032'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.06%  '
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '1.831.87'
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.260'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.372'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.63'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.25%  '
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008568'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").Value = '26.647.06'
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.023'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.913'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.90'
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.92'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.001'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.70'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.875'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08816'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.153'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.847'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.162'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7350'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.430'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.074'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.951'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05153'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.867'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5063'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.074'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.008'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4657'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.04'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.65%  '
$ws.Range("E47").Value = '  -1.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.572'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06023'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.97'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.83'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.57%  '
